$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price and 1h volume change), matching the
# refreshed feed from coinranking.com. A couple of rows were also
# re-ranked (swapped order) as part of the refresh.

$ws.Cells.Item(2, 4).Value = '27.691.63'
$ws.Cells.Item(2, 5).Value = '  -0.88%  '

$ws.Cells.Item(3, 4).Value = '1.623.02'
$ws.Cells.Item(3, 5).Value = '  -1.36%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.992'
$ws.Cells.Item(4, 5).Value = '  -0.82%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '210.73'
$ws.Cells.Item(5, 5).Value = '  -1.28%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.520'
$ws.Cells.Item(6, 5).Value = '  -1.18%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.991'
$ws.Cells.Item(7, 5).Value = '  -0.85%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '23.15'
$ws.Cells.Item(8, 5).Value = '  -1.19%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.257'
$ws.Cells.Item(9, 5).Value = '  -3.08%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0611'
$ws.Cells.Item(10, 5).Value = '  -0.59%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0876'
$ws.Cells.Item(11, 5).Value = '  +0.34%  '

$ws.Cells.Item(12, 4).Value = '1.858.34'
$ws.Cells.Item(12, 5).Value = '  -1.10%  '

$ws.Cells.Item(13, 4).Value = '1.625.32'
$ws.Cells.Item(13, 5).Value = '  -1.18%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.03'
$ws.Cells.Item(14, 5).Value = '  -1.03%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.561'
$ws.Cells.Item(15, 5).Value = '  -0.58%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '65.01'
$ws.Cells.Item(16, 5).Value = '  -0.73%  '

$ws.Cells.Item(17, 4).Value = '27.736.90'
$ws.Cells.Item(17, 5).Value = '  -0.74%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '229.71'
$ws.Cells.Item(18, 5).Value = '  -0.60%  '

$ws.Cells.Item(19, 5).Value = '  -0.54%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.48'
$ws.Cells.Item(20, 5).Value = '  -2.21%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.992'
$ws.Cells.Item(21, 5).Value = '  -0.80%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '10.33'
$ws.Cells.Item(22, 5).Value = '  -3.05%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '4.34'
$ws.Cells.Item(23, 5).Value = '  -1.30%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.04'
$ws.Cells.Item(24, 5).Value = '  -4.51%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '154.10'
$ws.Cells.Item(25, 5).Value = '  +1.35%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '6.92'
$ws.Cells.Item(26, 5).Value = '  +0.07%  '

$ws.Cells.Item(27, 2).Value = 'Stellar'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.111'
$ws.Cells.Item(27, 5).Value = '  -1.20%  '

$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '15.57'
$ws.Cells.Item(28, 5).Value = '  -1.01%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.993'
$ws.Cells.Item(29, 5).Value = '  -0.80%  '

$ws.Cells.Item(30, 5).Value = '  -1.58%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.0480'
$ws.Cells.Item(31, 5).Value = '  -1.07%  '

$ws.Cells.Item(32, 5).Value = '  +2.26%  '

$ws.Cells.Item(33, 4).Value = '1.398.29'
$ws.Cells.Item(33, 5).Value = '  -2.92%  '

$ws.Cells.Item(34, 5).Value = '  +0.06%  '

$ws.Cells.Item(35, 5).Value = '  -0.23%  '

$ws.Cells.Item(36, 5).Value = '  +7.36%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.34'
$ws.Cells.Item(37, 5).Value = '  +0.63%  '

$ws.Cells.Item(38, 5).Value = '  +0.27%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.556'
$ws.Cells.Item(39, 5).Value = '  -0.67%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.863'
$ws.Cells.Item(40, 5).Value = '  -3.02%  '

$ws.Cells.Item(41, 5).Value = '  -0.47%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.992'
$ws.Cells.Item(42, 5).Value = '  -0.74%  '

$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '66.35'
$ws.Cells.Item(43, 5).Value = '  -3.92%  '

$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '5.49'
$ws.Cells.Item(44, 5).Value = '  +1.24%  '

$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '1.82'
$ws.Cells.Item(45, 5).Value = '  -0.19%  '

$ws.Cells.Item(46, 2).Value = 'MXToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.19'
$ws.Cells.Item(46, 5).Value = '  -1.35%  '

$ws.Cells.Item(47, 4).Value = '1.766.57'
$ws.Cells.Item(47, 5).Value = '  -1.21%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '87.75'
$ws.Cells.Item(48, 5).Value = '  -1.54%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.0994'
$ws.Cells.Item(49, 5).Value = '  -1.47%  '

$ws.Cells.Item(50, 5).Value = '  -0.42%  '

$ws.Cells.Item(51, 4).Value = '0.0₇0967'
$ws.Cells.Item(51, 5).Value = '  -9.02%  '
